# Implementation revisions to budget by model number and modify model
# output written to Excel (i.e. oh, rates, asset info, formulas vs hard code)

$wb = $excel.ActiveWorkbook

# ----------------------------------------------------------------------
# Sheet "BME_DI_BCCW" (second tab)
# ----------------------------------------------------------------------
$wsW = $wb.Worksheets.Item("BME_DI_BCCW")
$wsW.Activate()

# actual_parts_exp (E) picks up what used to live in actual_parts... (F),
# and actual_parts_exp's old hard-coded "0" becomes the prior F value.
# budgeted_partial_oh / actual_partial_oh formulas now subtract the parts
# column too.
$wsW.Range("B2").Formula  = "=C2-D2-E2"
$wsW.Range("B3:B6").Formula = "=C3-D3-E3"
$wsW.Range("G2").Formula  = "=H2-I2-J2"
$wsW.Range("G3:G6").Formula = "=H3-I3-J3"

$wsW.Range("E2").Value = 146756.97
$wsW.Range("E3").Value = 106128.98
$wsW.Range("E4").Value = 491988.52
$wsW.Range("E5").Value = 837300.31
$wsW.Range("E6").Value = 1267804.3600000001

$wsW.Range("F2").Value = 426564.27
$wsW.Range("F3").Value = 224633.75
$wsW.Range("F4").Value = 73137.38
$wsW.Range("F5").Value = 30713.14
$wsW.Range("F6").Value = 174491.42

# F2's style loses its border (s=5 -> s=7) to match the other "no border"
# currency cells (e.g. H5).
$wsW.Range("F2").Borders.LineStyle = -4142

$wsW.Range("J2").Value = 246732
$wsW.Range("J3").Value = 246732
$wsW.Range("J4").Value = 246732
$wsW.Range("J5").Value = 500000
$wsW.Range("J6").Value = 246732

$wsW.Range("K2").Value = 341000
$wsW.Range("K3").Value = 341000
$wsW.Range("K4").Value = 341000
$wsW.Range("K5").Value = 87732
$wsW.Range("K6").Value = 341000

$wsW.Range("H12").Select()

# ----------------------------------------------------------------------
# Sheet "BME_DI_BCC" (third tab)
# ----------------------------------------------------------------------
$wsC = $wb.Worksheets.Item("BME_DI_BCC")
$wsC.Activate()

$wsC.Range("B2").Formula  = "=C2-D2-E2"
$wsC.Range("B3:B6").Formula = "=C3-D3-E3"
$wsC.Range("G2").Formula  = "=H2-I2-J2"
$wsC.Range("G3:G6").Formula = "=H3-I3-J3"

$wsC.Range("E2").Value = 57962.559999999998
$wsC.Range("E3").Value = 94021.16
$wsC.Range("E4").Value = 80898.429999999993
$wsC.Range("E5").Value = 248529.01
$wsC.Range("E6").Value = 203994.44

$wsC.Range("F2").Value = 128509.74
$wsC.Range("F3").Value = 107835.51
$wsC.Range("F4").Value = 82756.039999999994
$wsC.Range("F5").Value = 75668.179999999993
$wsC.Range("F6").Value = 44033.97

$wsC.Range("J2").Value = 60000
$wsC.Range("J3").Value = 60000
$wsC.Range("J4").Value = 60000
$wsC.Range("J5").Value = 80000
$wsC.Range("J6").Value = 165000

$wsC.Range("K2").Value = 217451
$wsC.Range("K3").Value = 217451
$wsC.Range("K4").Value = 217451
$wsC.Range("K5").Value = 197451
$wsC.Range("K6").Value = 112451

$wsC.Range("D36").Select()

# BME_DI_BCC ends up the active tab / active sheet.
$wsC.Activate()
